# ----------------------------------------------------------------------------
# Applies the commit "Fixed update to excel issue":
#   1. Renames the "Requested quantity" header on "Weekly Quantity"  -> "Weekly_PO_Qty"
#   2. Renames the "Requested quantity" header on "Monthly Trend"    -> "Monthly_PO_Qty"
#   3. Adds a new "PO Forecast" worksheet (sheetId 3, placed after "Monthly Trend")
#      with headers ds / PO_Forecast / yhat_lower / yhat_upper and 54 data rows.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$weekly  = $wb.Worksheets.Item("Weekly Quantity")
$monthly = $wb.Worksheets.Item("Monthly Trend")

# 1. + 2. Rename the header cells in-place (keep existing bold/bordered style).
$weekly.Range("B1").Value  = "Weekly_PO_Qty"
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# 3. Add the new "PO Forecast" sheet right after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "PO Forecast"

# Match the page margins used by the other sheets (values are in points;
# 72pt = 1in, so 54/72/36pt give the 0.75in/1in/0.5in margins seen elsewhere).
$ws.PageSetup.LeftMargin   = 54
$ws.PageSetup.RightMargin  = 54
$ws.PageSetup.TopMargin    = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Reuse the exact header style (bold + border + centered) from column A of an
# existing sheet, broadcast across the whole header row, then overwrite text.
$weekly.Range("A1").Copy($ws.Range("A1:D1"))
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# Reuse the exact date/time number-format style used for column A's data
# cells, broadcast down the whole date column.
$weekly.Range("A2").Copy($ws.Range("A2:A55"))

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$ws.Range("A2").Value = 44934.99999999999
$ws.Range("B2").Value = 13
$ws.Range("C2").Value = -0.5271067943040574
$ws.Range("D2").Value = 25.40147441483534
$ws.Range("A3").Value = 44941.99999999999
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = -0.1820261687052689
$ws.Range("D3").Value = 26.74968295428741
$ws.Range("A4").Value = 44955.99999999999
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = -0.8198697717672082
$ws.Range("D4").Value = 25.41534075317099
$ws.Range("A5").Value = 44962.99999999999
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = -0.820406925766097
$ws.Range("D5").Value = 24.67688120580234
$ws.Range("A6").Value = 44969.99999999999
$ws.Range("B6").Value = 13
$ws.Range("C6").Value = -1.528659376095897
$ws.Range("D6").Value = 26.24007510919779
$ws.Range("A7").Value = 44976.99999999999
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = -0.6020844320315634
$ws.Range("D7").Value = 25.95549432301565
$ws.Range("A8").Value = 44983.99999999999
$ws.Range("B8").Value = 13
$ws.Range("C8").Value = 0.637522609016531
$ws.Range("D8").Value = 24.99781460982209
$ws.Range("A9").Value = 45011.99999999999
$ws.Range("B9").Value = 13
$ws.Range("C9").Value = 0.5296443744195011
$ws.Range("D9").Value = 26.91750307488106
$ws.Range("A10").Value = 45032.99999999999
$ws.Range("B10").Value = 13
$ws.Range("C10").Value = -0.3528492402394646
$ws.Range("D10").Value = 25.5541118333717
$ws.Range("A11").Value = 45039.99999999999
$ws.Range("B11").Value = 13
$ws.Range("C11").Value = 0.350886568724831
$ws.Range("D11").Value = 25.42678843105342
$ws.Range("A12").Value = 45060.99999999999
$ws.Range("B12").Value = 12
$ws.Range("C12").Value = -0.7293934447428303
$ws.Range("D12").Value = 25.17472578662952
$ws.Range("A13").Value = 45067.99999999999
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 0.4046227444137188
$ws.Range("D13").Value = 24.75766744706643
$ws.Range("A14").Value = 45074.99999999999
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 0.4640824661565389
$ws.Range("D14").Value = 25.52246460860911
$ws.Range("A15").Value = 45088.99999999999
$ws.Range("B15").Value = 12
$ws.Range("C15").Value = -1.663573350724932
$ws.Range("D15").Value = 25.14371805180408
$ws.Range("A16").Value = 45095.99999999999
$ws.Range("B16").Value = 12
$ws.Range("C16").Value = -1.803894246172023
$ws.Range("D16").Value = 26.08037310398392
$ws.Range("A17").Value = 45109.99999999999
$ws.Range("B17").Value = 12
$ws.Range("C17").Value = -1.16063972480592
$ws.Range("D17").Value = 26.00824043806533
$ws.Range("A18").Value = 45116.99999999999
$ws.Range("B18").Value = 12
$ws.Range("C18").Value = -1.505997954489133
$ws.Range("D18").Value = 25.50710402975074
$ws.Range("A19").Value = 45123.99999999999
$ws.Range("B19").Value = 12
$ws.Range("C19").Value = -0.4760166180176912
$ws.Range("D19").Value = 25.02413972559179
$ws.Range("A20").Value = 45130.99999999999
$ws.Range("B20").Value = 12
$ws.Range("C20").Value = 0.03406800192288839
$ws.Range("D20").Value = 25.90339339010362
$ws.Range("A21").Value = 45137.99999999999
$ws.Range("B21").Value = 12
$ws.Range("C21").Value = -1.163823939970937
$ws.Range("D21").Value = 25.32517241857733
$ws.Range("A22").Value = 45151.99999999999
$ws.Range("B22").Value = 12
$ws.Range("C22").Value = -0.2825779226418432
$ws.Range("D22").Value = 26.57453247841655
$ws.Range("A23").Value = 45179.99999999999
$ws.Range("B23").Value = 12
$ws.Range("C23").Value = -0.2046342919008155
$ws.Range("D23").Value = 25.58803138475572
$ws.Range("A24").Value = 45186.99999999999
$ws.Range("B24").Value = 12
$ws.Range("C24").Value = -0.4807511891047375
$ws.Range("D24").Value = 25.43667148080623
$ws.Range("A25").Value = 45193.99999999999
$ws.Range("B25").Value = 12
$ws.Range("C25").Value = 0.4885821710578017
$ws.Range("D25").Value = 26.48037604382117
$ws.Range("A26").Value = 45200.99999999999
$ws.Range("B26").Value = 12
$ws.Range("C26").Value = -0.2361148844738142
$ws.Range("D26").Value = 26.85756042061587
$ws.Range("A27").Value = 45207.99999999999
$ws.Range("B27").Value = 12
$ws.Range("C27").Value = -1.047336028127639
$ws.Range("D27").Value = 26.32307868723293
$ws.Range("A28").Value = 45214.99999999999
$ws.Range("B28").Value = 12
$ws.Range("C28").Value = -0.2408813597711066
$ws.Range("D28").Value = 24.11214652606886
$ws.Range("A29").Value = 45221.99999999999
$ws.Range("B29").Value = 12
$ws.Range("C29").Value = -0.1857030724138984
$ws.Range("D29").Value = 25.74620181319534
$ws.Range("A30").Value = 45235.99999999999
$ws.Range("B30").Value = 12
$ws.Range("C30").Value = -1.979548033969344
$ws.Range("D30").Value = 25.2873130389532
$ws.Range("A31").Value = 45242.99999999999
$ws.Range("B31").Value = 12
$ws.Range("C31").Value = -1.028417104254006
$ws.Range("D31").Value = 24.54903512591317
$ws.Range("A32").Value = 45249.99999999999
$ws.Range("B32").Value = 12
$ws.Range("C32").Value = -1.232863451105629
$ws.Range("D32").Value = 25.22649682486739
$ws.Range("A33").Value = 45256.99999999999
$ws.Range("B33").Value = 12
$ws.Range("C33").Value = -0.8534862539136431
$ws.Range("D33").Value = 25.13990741130534
$ws.Range("A34").Value = 45270.99999999999
$ws.Range("B34").Value = 12
$ws.Range("C34").Value = -0.5869631302885987
$ws.Range("D34").Value = 25.72788331135697
$ws.Range("A35").Value = 45298.99999999999
$ws.Range("B35").Value = 12
$ws.Range("C35").Value = -0.5866891250688547
$ws.Range("D35").Value = 25.70736196486773
$ws.Range("A36").Value = 45305.99999999999
$ws.Range("B36").Value = 12
$ws.Range("C36").Value = 0.1150482964867209
$ws.Range("D36").Value = 24.96412126799478
$ws.Range("A37").Value = 45312.99999999999
$ws.Range("B37").Value = 12
$ws.Range("C37").Value = -0.8789801709285051
$ws.Range("D37").Value = 24.73781481882762
$ws.Range("A38").Value = 45319.99999999999
$ws.Range("B38").Value = 12
$ws.Range("C38").Value = -0.947544736075804
$ws.Range("D38").Value = 25.20817045242906
$ws.Range("A39").Value = 45333.99999999999
$ws.Range("B39").Value = 12
$ws.Range("C39").Value = -1.105138571388012
$ws.Range("D39").Value = 25.145434455711
$ws.Range("A40").Value = 45382.99999999999
$ws.Range("B40").Value = 12
$ws.Range("C40").Value = -0.3194105872054248
$ws.Range("D40").Value = 25.02268569489624
$ws.Range("A41").Value = 45389.99999999999
$ws.Range("B41").Value = 12
$ws.Range("C41").Value = -1.142578754124113
$ws.Range("D41").Value = 24.77594716038065
$ws.Range("A42").Value = 45396.99999999999
$ws.Range("B42").Value = 12
$ws.Range("C42").Value = -0.6558614891311449
$ws.Range("D42").Value = 25.25873085631562
$ws.Range("A43").Value = 45403.99999999999
$ws.Range("B43").Value = 12
$ws.Range("C43").Value = -0.2246914165477478
$ws.Range("D43").Value = 25.50784610143367
$ws.Range("A44").Value = 45459.99999999999
$ws.Range("B44").Value = 12
$ws.Range("C44").Value = -0.6895275934550787
$ws.Range("D44").Value = 25.69307790768651
$ws.Range("A45").Value = 45487.99999999999
$ws.Range("B45").Value = 12
$ws.Range("C45").Value = -0.347690129828965
$ws.Range("D45").Value = 25.32699970853373
$ws.Range("A46").Value = 45515.99999999999
$ws.Range("B46").Value = 12
$ws.Range("C46").Value = -1.74322569508891
$ws.Range("D46").Value = 26.45934714805481
$ws.Range("A47").Value = 45543.99999999999
$ws.Range("B47").Value = 12
$ws.Range("C47").Value = -0.9148155663036598
$ws.Range("D47").Value = 24.58926914803034
$ws.Range("A48").Value = 45550.99999999999
$ws.Range("B48").Value = 12
$ws.Range("C48").Value = -0.5852531395234977
$ws.Range("D48").Value = 25.04604235682412
$ws.Range("A49").Value = 45557.99999999999
$ws.Range("B49").Value = 12
$ws.Range("C49").Value = -0.5138890992103153
$ws.Range("D49").Value = 25.57852185352675
$ws.Range("A50").Value = 45564.99999999999
$ws.Range("B50").Value = 12
$ws.Range("C50").Value = 0.1350148653226725
$ws.Range("D50").Value = 25.20575479775067
$ws.Range("A51").Value = 45571.99999999999
$ws.Range("B51").Value = 12
$ws.Range("C51").Value = -0.09565960846064499
$ws.Range("D51").Value = 25.39295363262861
$ws.Range("A52").Value = 45578.99999999999
$ws.Range("B52").Value = 12
$ws.Range("C52").Value = -1.01999409112488
$ws.Range("D52").Value = 24.39551406142128
$ws.Range("A53").Value = 45585.99999999999
$ws.Range("B53").Value = 12
$ws.Range("C53").Value = -0.466181971532182
$ws.Range("D53").Value = 25.02646936159381
$ws.Range("A54").Value = 45592.99999999999
$ws.Range("B54").Value = 12
$ws.Range("C54").Value = -0.4054136803213905
$ws.Range("D54").Value = 24.85715305091874
$ws.Range("A55").Value = 45599.99999999999
$ws.Range("B55").Value = 12
$ws.Range("C55").Value = -0.7390450993668164
$ws.Range("D55").Value = 26.30186801033419

# Restore the originally-active sheet (adding a sheet makes it the active one).
$weekly.Activate()
